# "spring v fall" — adds the spring-vs-fall zooplankton CPUE regression (lm)
# and PERMANOVA (adonis) report as a new worksheet, matching the Lucida-Console
# "R console output" styling already used on the "pairwise" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "pairwise" — source of the existing report styles

# Style source cells on "pairwise":
#   A9 -> blank divider row style
#   A8 -> plain Lucida-Console body-text style
#   A7 -> Lucida-Console "section total" style (white fill)
$styleBlank = $ws1.Range("A9")
$styleBody  = $ws1.Range("A8")
$styleTotal = $ws1.Range("A7")

# Add the new worksheet at the end of the tab strip and rename it
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "spring v fall"

$styleBody.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 'Call:'

$styleBody.Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 'lm(formula = log(tCPUE) ~ site + targets2 + season, data = bugsblitzSF.1[which(bugsblitzSF.1$SampleID != '

$styleBody.Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = '    "MAC6-22MAR2018"), ])'

$styleBlank.Copy()
$ws.Range("A5").PasteSpecial(-4122)

$styleBody.Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 'Residuals:'

$styleBody.Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = '    Min      1Q  Median      3Q     Max '

$styleBody.Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = '-4.3307 -1.3386 -0.1527  1.1829  4.0837 '

$styleBlank.Copy()
$ws.Range("A9").PasteSpecial(-4122)

$styleBody.Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 'Coefficients:'

$styleBody.Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = '                  Estimate Std. Error t value Pr(>|t|)  '

$styleBody.Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = '(Intercept)         1.4008     0.5835   2.401   0.0192 *'

$styleBody.Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 'siteBrowns         -1.5584     0.6896  -2.260   0.0272 *'

$styleBody.Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 'siteWinter         -1.2670     0.7046  -1.798   0.0768 .'

$styleBody.Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 'siteProspect        1.6347     0.6740   2.425   0.0181 *'

$styleBody.Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 'targets2sweep net   0.4973     0.5075   0.980   0.3307  '

$styleBody.Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 'seasonfall          1.0949     0.5002   2.189   0.0322 *'

$styleBody.Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = '---'

$styleBody.Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = 'Signif. codes:  0 ‘***’ 0.001 ‘**’ 0.01 ‘*’ 0.05 ‘.’ 0.1 ‘ ’ 1'

$styleBlank.Copy()
$ws.Range("A20").PasteSpecial(-4122)

$styleBody.Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = 'Residual standard error: 2.018 on 65 degrees of freedom'

$styleBody.Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = 'Multiple R-squared:  0.3436,    Adjusted R-squared:  0.2931 '

$styleTotal.Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = 'F-statistic: 6.805 on 5 and 65 DF,  p-value: 3.692e-05'

$styleBody.Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = 'adonis(formula = sf.12p ~ site + targets2 + season, data = sf2) '

$styleBlank.Copy()
$ws.Range("A27").PasteSpecial(-4122)

$styleBody.Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = 'Permutation: free'

$styleBody.Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 'Number of permutations: 999'

$styleBlank.Copy()
$ws.Range("A30").PasteSpecial(-4122)

$styleBody.Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = 'Terms added sequentially (first to last)'

$styleBlank.Copy()
$ws.Range("A32").PasteSpecial(-4122)

$styleBody.Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = '          Df SumsOfSqs MeanSqs F.Model      R2 Pr(>F)    '

$styleBody.Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A34").Value = 'site       3    4.6812 1.56039 10.0981 0.25666  0.001 ***'

$styleBody.Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = 'targets2   1    0.5781 0.57806  3.7409 0.03169  0.011 *  '

$styleBody.Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = 'season     1    2.7813 2.78127 17.9990 0.15249  0.001 ***'

$styleBody.Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = 'Residuals 66   10.1986 0.15452         0.55916           '

$styleTotal.Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A38").Value = 'Total     71   18.2391                 1.00000           '

$excel.CutCopyMode = $false

# Leave the new sheet active with the same selection Excel recorded on save
$ws.Range("J30").Select()
